# 977: Update xlsx data
# Update GS sheet OM_Key values and OM_Team_Key values, resize column H,
# and move the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GS")

# Row 2: OM_Key 1357 -> 1001, OM_Team_Key "T1" -> "WMT"
$ws.Cells.Item(2, 6).Value = 1001
$ws.Cells.Item(2, 8).Value = "WMT"

# Row 3: OM_Key 9876 -> 1002, OM_Team_Key "T1" -> "WMT"
$ws.Cells.Item(3, 6).Value = 1002
$ws.Cells.Item(3, 8).Value = "WMT"

# Widen column H (OM_Team_Key) to fit the new text
$ws.Columns.Item(8).ColumnWidth = 17

# Move the selected/active cell to H4
$ws.Range("H4").Select()
